$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 (the "hk02" plane entry). All rows below shift up by one,
# and the sheet dimension shrinks from A1:J9 to A1:J8.
$ws.Rows.Item(3).Delete()
